# Update the "想去人数" (F column) figures on both the "展览" and
# "全部类型" worksheets, which contain duplicate data tables.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# row -> new F value
$updates = @{
    2  = 8393
    3  = 7955
    9  = 130
    14 = 1988
    16 = 61
    20 = 27
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
